$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.616.40"
$ws.Range("E2").Value = "  -4.49%  "

$ws.Range("D3").Value = "2.936.01"
$ws.Range("E3").Value = "  -2.36%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'549.25"
$ws.Range("E5").Value = "  -4.51%  "

$ws.Range("D6").Value = "'130.38"
$ws.Range("E6").Value = "  +3.24%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D9").Value = "2.929.55"
$ws.Range("E9").Value = "  -2.44%  "

$ws.Range("D10").Value = "'0.126"
$ws.Range("E10").Value = "  -4.42%  "

$ws.Range("D11").Value = "'4.76"
$ws.Range("E11").Value = "  -6.10%  "

$ws.Range("E12").Value = "  +1.04%  "

$ws.Range("E13").Value = "  -0.47%  "

$ws.Range("D14").Value = "'32.86"
$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("E15").Value = "  +0.41%  "

$ws.Range("D16").Value = "3.420.18"
$ws.Range("E16").Value = "  -2.33%  "

$ws.Range("E17").Value = "  +6.22%  "

$ws.Range("D18").Value = "2.929.14"

$ws.Range("D19").Value = "57.639.10"
$ws.Range("E19").Value = "  -4.26%  "

$ws.Range("D20").Value = "'416.80"
$ws.Range("E20").Value = "  -3.08%  "

$ws.Range("D21").Value = "'13.17"
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").Value = "'0.687"
$ws.Range("E22").Value = "  +2.53%  "

$ws.Range("E23").Value = "  -1.39%  "

$ws.Range("D24").Value = "'13.05"
$ws.Range("E24").Value = "  +0.85%  "

$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("E28").Value = "  -3.11%  "

$ws.Range("D29").Value = "'7.49"
$ws.Range("E29").Value = "  +2.65%  "

$ws.Range("E30").Value = "  +1.14%  "

$ws.Range("E31").Value = "  -0.85%  "

$ws.Range("D32").Value = "'5.99"
$ws.Range("E32").Value = "  -2.78%  "

$ws.Range("D33").Value = "'0.0967"
$ws.Range("E33").Value = "  +2.28%  "

$ws.Range("E34").Value = "  +0.76%  "

$ws.Range("D35").Value = "'0.941"
$ws.Range("E35").Value = "  +0.13%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").Value = "'48.14"
$ws.Range("E37").Value = "  -4.44%  "

$ws.Range("D38").Value = "0.0₃0682"
$ws.Range("E38").Value = "  +0.83%  "

$ws.Range("E39").Value = "  +2.58%  "

$ws.Range("E40").Value = "  +3.00%  "

$ws.Range("D41").Value = "'379.10"
$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("D42").Value = "'0.108"
$ws.Range("E42").Value = "  -0.81%  "

$ws.Range("E43").Value = "  -3.23%  "

$ws.Range("D44").Value = "2.687.61"
$ws.Range("E44").Value = "  +0.36%  "

$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("E46").Value = "  +1.19%  "

$ws.Range("D47").Value = "'121.95"
$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("E48").Value = "  +1.39%  "

$ws.Range("D49").Value = "'1.97"
$ws.Range("E49").Value = "  -1.96%  "

$ws.Range("D50").Value = "'23.07"
$ws.Range("E50").Value = "  -2.33%  "

$ws.Range("E51").Value = "  -0.61%  "
